# Update Back/Lay odds (columns F:AO) for rows 2-12 on Sheet1
# to match the new odds snapshot for 2026-01-02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 11,36
$row0 = @(5.4,5.6,2.14,2.16,2.8,2.84,0,1.25,1.8,2.2,1.24,4.9,1.06,15,3.45,1.36,1.84,1.21,4.9,4.7,9.6,34,10,8.199999999999999,17.5,75,40,36,85,280,260,260,540,1000,1000,85)
for ($j = 0; $j -lt 36; $j++) { $data[0,$j] = $row0[$j] }
$row1 = @(1.71,1.72,6.2,6.4,3.85,3.95,1.45,1.09,3.35,1.41,1.77,2.22,1.28,4.3,1.05,1.01,1.2,2.28,12.5,18.5,55,220,7,8.6,24,120,9.199999999999999,10.5,25,120,17,21,46,210,14,190)
for ($j = 0; $j -lt 36; $j++) { $data[1,$j] = $row1[$j] }
$row2 = @(2.98,3.2,2.72,2.8,3.1,3.25,1.76,1.17,2.2,1.74,1.37,3.35,1.12,8.199999999999999,2.26,1.04,1.55,1.46,12.5,7.6,16,190,1000,42,16.5,230,1000,1000,34,120,1000,70,140,1000,110,85)
for ($j = 0; $j -lt 36; $j++) { $data[2,$j] = $row2[$j] }
$row3 = @(2.36,2.48,3.4,3.75,3.2,3.4,1.48,1.1,3.3,1.42,1.77,2.18,1.29,4.1,1.86,1.99,1.37,1.67,12,13,26,250,9.4,7.6,15,48,80,12,18.5,330,85,46,120,1000,25,75)
for ($j = 0; $j -lt 36; $j++) { $data[3,$j] = $row3[$j] }
$row4 = @(4.1,4.3,1.79,1.82,4.6,4.8,1.23,1.02,8.199999999999999,1.12,3.4,1.39,1.97,1.99,1.46,3,2.2,1.3,42,18.5,17.5,23,30,29,14,44,42,18,14.5,22,130,42,34,46,1000,5.9)
for ($j = 0; $j -lt 36; $j++) { $data[4,$j] = $row4[$j] }
$row5 = @(8,8.199999999999999,1.52,1.53,4.5,4.6,1.44,1.07,3.6,1.36,1.87,2.1,1.32,3.85,2.22,1.77,2.88,1.13,15.5,6.8,7.6,13,22,9.800000000000001,10,17,70,32,29,50,310,150,140,200,260,9.800000000000001)
for ($j = 0; $j -lt 36; $j++) { $data[5,$j] = $row5[$j] }
$row6 = @(1.58,1.61,7.6,8.4,4,4.2,1.54,1.1,3.1,1.45,1.71,2.36,1.26,4.6,2.38,1.68,1.13,2.64,11,19,70,320,6.4,9.199999999999999,32,180,7.8,10.5,32,230,14.5,20,60,270,13,1000)
for ($j = 0; $j -lt 36; $j++) { $data[6,$j] = $row6[$j] }
$row7 = @(2.7,2.72,2.9,2.94,3.45,3.5,1.45,1.08,3.75,1.35,1.91,2.06,1.35,3.7,1.8,2.2,1.51,1.58,13,11.5,18,46,11,7.4,12.5,32,17,12,17,46,40,29,40,90,27,30)
for ($j = 0; $j -lt 36; $j++) { $data[7,$j] = $row7[$j] }
$row8 = @(7,7.2,1.63,1.64,4.1,4.2,1.48,1.08,3.4,1.4,1.81,2.2,1.31,4.1,2.2,1.8,2.58,1.16,12.5,6.8,8.199999999999999,15,19,9,10,18,55,26,25,46,230,120,120,170,180,11.5)
for ($j = 0; $j -lt 36; $j++) { $data[8,$j] = $row8[$j] }
$row9 = @(2.12,2.14,5,5.1,3,3.05,1.78,1.18,2.22,1.8,1.39,3.45,1.12,8.199999999999999,2.72,1.55,1.24,1.87,6.4,10.5,36,160,5.4,7.4,24,120,10,13,34,170,27,36,95,370,38,230)
for ($j = 0; $j -lt 36; $j++) { $data[9,$j] = $row9[$j] }
$row10 = @(1.95,1.96,4.6,4.7,3.6,3.7,1.48,1.09,3.35,1.41,1.79,2.22,1.29,4.2,1.98,1.96,1.27,2.04,12,14.5,34,120,7.8,8,19,70,11,10.5,22,80,24,23,42,140,17,85)
for ($j = 0; $j -lt 36; $j++) { $data[10,$j] = $row10[$j] }

$ws.Range("F2:AO12").Value = $data
